$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.589.67"
$ws.Range("E2").Value = "  +2.08%  "
$ws.Range("D3").Value = "1.688.12"
$ws.Range("E3").Value = "  +3.34%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "220.74"
$ws.Range("E5").Value = "  +2.64%  "
$ws.Range("E6").Value = "  +0.35%  "
$ws.Range("E7").Value = "  -0.26%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "30.95"
$ws.Range("E8").Value = "  +4.36%  "
$ws.Range("E9").Value = "  +2.34%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0626"
$ws.Range("E10").Value = "  +2.05%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0903"
$ws.Range("E11").Value = "  -1.52%  "
$ws.Range("E12").Value = "  +3.45%  "
$ws.Range("E13").Value = "  +13.83%  "
$ws.Range("E14").Value = "  +8.12%  "
$ws.Range("D15").Value = "1.680.00"
$ws.Range("E15").Value = "  +3.17%  "
$ws.Range("E16").Value = "  +3.24%  "
$ws.Range("D17").Value = "30.603.15"
$ws.Range("E17").Value = "  +2.10%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "66.08"
$ws.Range("E18").Value = "  +1.80%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "248.76"
$ws.Range("E19").Value = "  -0.25%  "
$ws.Range("E20").Value = "  +1.69%  "
$ws.Range("E21").Value = "  -0.25%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.32"
$ws.Range("E22").Value = "  +3.46%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.23"
$ws.Range("E23").Value = "  +5.66%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.20"
$ws.Range("E24").Value = "  +3.51%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "157.32"
$ws.Range("E25").Value = "  -1.71%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.98"
$ws.Range("E26").Value = "  +1.56%  "
$ws.Range("E27").Value = "  +0.32%  "
$ws.Range("E28").Value = "  +2.20%  "
$ws.Range("E29").Value = "  -0.24%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0501"
$ws.Range("E30").Value = "  +2.00%  "
$ws.Range("E31").Value = "  +1.01%  "
$ws.Range("E32").Value = "  +3.02%  "
$ws.Range("E33").Value = "  +2.68%  "
$ws.Range("D34").Value = "1.510.46"
$ws.Range("E34").Value = "  +5.33%  "
$ws.Range("E35").Value = "  +5.03%  "
$ws.Range("E36").Value = "  -0.51%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0180"
$ws.Range("E37").Value = "  +4.43%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "79.75"
$ws.Range("E38").Value = "  +8.29%  "
$ws.Range("B39").Value = "ImmutableX"
$ws.Range("C39").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.586"
$ws.Range("E39").Value = "  +5.20%  "
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.71"
$ws.Range("E40").Value = "  -5.47%  "
$ws.Range("E41").Value = "  +1.33%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.855"
$ws.Range("E42").Value = "  +2.71%  "
$ws.Range("E43").Value = "  +1.11%  "
$ws.Range("E44").Value = "  +1.57%  "
$ws.Range("E45").Value = "  -2.10%  "
$ws.Range("E46").Value = "  -0.28%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "52.49"
$ws.Range("E47").Value = "  -4.36%  "
$ws.Range("E48").Value = "  +2.81%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.44"
$ws.Range("E49").Value = "  -0.83%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "95.67"
$ws.Range("E50").Value = "  +5.79%  "
$ws.Range("E51").Value = "  +7.23%  "
